$d = $word.ActiveDocument

# 1) Insert a new bold paragraph "Showing changes between current and last."
#    right after the first (empty) paragraph of the document.
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphAfter()
$newHeaderPara = $d.Paragraphs(2)
$newHeaderPara.Range.Text = "Showing changes between current and last."
$headerRange = $newHeaderPara.Range
$headerTextRange = $d.Range($headerRange.Start, $headerRange.End - 1)
$headerTextRange.Font.Bold = 1

# 2) Insert a new paragraph "Mars Missing" right before the "SoT Missing"
#    paragraph that sits under "Overall Resource DCA - Now Confidence changes".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^SoT Missing\r?$") {
        $target = $p
        break
    }
}
$target.Range.InsertParagraphBefore()
# After InsertParagraphBefore, $target now refers to the freshly-created
# (empty) paragraph that was inserted immediately before the original
# "SoT Missing" paragraph, so we set its text directly.
$target.Range.Text = "Mars Missing"

# 3) Bump the missing-rating total for that same section from 5 to 6.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^5 project\(s\) in total are missing a rating") {
        $p.Range.Text = "6 project(s) in total are missing a rating"
    }
}
